$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: item 7 "FAKTUR" KIRIM line, diff against row 22/18 total ---
$ws.Range("A26").Value = 7
$ws.Range("C26").Value = "FAKTUR"
$ws.Range("E26").Value = 410
$ws.Range("F26").Formula = "=F22-E26"

# --- Row 27 (replaces the old stray F27 total) ---
$ws.Range("E27").Value = 410
$ws.Range("F27").Formula = "=F23-E27"

# --- Row 28 (new) ---
$ws.Range("E28").Value = 410
$ws.Range("F28").Formula = "=F24-E28"

# --- Row 30: item 8 "6 MEI" kirim 70 set ---
$ws.Range("A30").Value = 8
$ws.Range("B30").Value = "6 MEI"
$ws.Range("C30").Value = "WEARPACK"
$ws.Range("D30").Value = 70
$ws.Range("F30").Formula = "=F26+D30"

# --- Row 31 (new) ---
$ws.Range("C31").Value = "KAOS"
$ws.Range("D31").Value = 70
$ws.Range("F31").Formula = "=F27+D31"

# --- Row 32 (new) ---
$ws.Range("C32").Value = "TOPI"
$ws.Range("D32").Value = 70
$ws.Range("F32").Formula = "=F28+D32"

# --- Update the view: zoom + selection per diff ---
$ws.Application.ActiveWindow.Zoom = 145
$null = $ws.Range("F10").Select()
